$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Query" column (D). This shifts the
# old D column (header "Query" + its formulas, along with their style) to E,
# and leaves a blank column D (inheriting column D's original formatting)
# ready for the new "PrimaryAddrID" data.
$ws.Range("D1").EntireColumn.Insert()

# New column header + per-row constant value (1) for PrimaryAddrID.
$ws.Range("D1").Value = "PrimaryAddrID"

$lastRow = 11
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("D$r").Value = 1

    # Rebuild the Query formula (now in column E) so it also pulls in the
    # new PrimaryAddrID column/header.
    $formula = '="INSERT INTO "&A' + $r + '&" ([" &B$1 &"],["&C$1&"],["&D$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'',''" & D' + $r + ' & "'')"'
    $ws.Range("E$r").Formula = $formula
}

$ws.Range("G11").Select()
